$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.281.32"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.511.44"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.32"
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.96"
$ws.Range("E6").Value = "  -1.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.510.96"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("E10").Value = "  -4.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.98"
$ws.Range("E11").Value = "  +2.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.412"
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.103.29"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000208"
$ws.Range("E14").Value = "  -2.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.35"
$ws.Range("E15").Value = "  -3.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.512.56"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.283.63"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.41"
$ws.Range("E19").Value = "  -1.82%  "
$ws.Range("E20").Value = "  -3.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.87"
$ws.Range("E21").Value = "  -3.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "426.15"
$ws.Range("E22").Value = "  -1.83%  "
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.01"
$ws.Range("E24").Value = "  -2.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.642.18"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.29"
$ws.Range("E28").Value = "  -5.19%  "
$ws.Range("E29").Value = "  -3.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.47"
$ws.Range("E30").Value = "  -1.06%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("E33").Value = "  -7.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.21"
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.499.83"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("E37").Value = "  -3.67%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.79"
$ws.Range("E38").Value = "  -2.71%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.61"
$ws.Range("E39").Value = "  -4.75%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "170.58"
$ws.Range("E41").Value = "  +1.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0859"
$ws.Range("E42").Value = "  -3.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.18"
$ws.Range("E43").Value = "  -4.38%  "
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("E45").Value = "  -8.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.39"
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("E47").Value = "  -8.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.84"
$ws.Range("E48").Value = "  -11.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.44"
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("E50").Value = "  -3.87%  "
$ws.Range("E51").Value = "  -3.02%  "
